# Reasigna los índices de punteo en orden secuencial, sin alterar el orden de las filas.
# La columna "Indice_Punteo" (A) se retira de su posición inicial y se reubica al final
# (columna Q), añadiendo además una nueva columna de índice secuencial (columna R).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Quitar la columna A ("Indice_Punteo"); el resto de columnas se desplaza una posición
#    a la izquierda (B->A, C->B, ..., Q->P), arrastrando consigo su formato.
$ws.Range("A1").EntireColumn.Delete()

# 2) Añadir las dos nuevas columnas de índice al final (Q y R), copiando el formato de
#    cabecera (negrita, bordes, alineación) de la última columna de datos (P1).
$ws.Range("P1").Copy()
$ws.Range("Q1:R1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("Q1").Value = "Indice_Punteo"
$ws.Range("R1").Value = "Indice_Punteo"

# Las celdas de datos para las nuevas columnas quedan vacías, pero presentes.
$ws.Range("Q2").Style = "Normal"
$ws.Range("R2").Style = "Normal"
